# Weekly update: insert a new week's worth of price rows (3 quality grades)
# for "Comercializadora del Agro de Limarí - Frutilla" at the top of the
# date-ordered block (rows 211-213), pushing the existing rows down by 3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows at 211 (existing rows 211.. shift down to 214..)
$ws.Range("A211:T213").EntireRow.Insert()

# Common column values shared by the whole "Comercializadora del Agro de
# Limarí" / "Frutilla" block.
$mercadoId = 2
$mercado   = "Comercializadora del Agro de Limarí"
$region    = "Coquimbo"
$fecha     = 44476
$codreg    = 4
$tipo      = "Fruta"
$productoId = 100101
$producto   = "Berries"
$categoriaId = 100112025
$categoria   = "Frutilla"
$variedad    = "Sin especificar"
$unidad      = "$/bandeja 7 kilos"
$origen      = "Provincia de Melipilla"
$kgUnidad    = 7

function Set-FrutillaRow($row, $calidad, $volumen, $precioMin, $precioMax, $precioProm, $precioKg) {
    $ws.Cells.Item($row, 1).Value = $mercadoId
    $ws.Cells.Item($row, 2).Value = $mercado
    $ws.Cells.Item($row, 3).Value = $region
    $ws.Cells.Item($row, 4).Value = $fecha
    $ws.Cells.Item($row, 5).Value = $codreg
    $ws.Cells.Item($row, 6).Value = $tipo
    $ws.Cells.Item($row, 7).Value = $productoId
    $ws.Cells.Item($row, 8).Value = $producto
    $ws.Cells.Item($row, 9).Value = $categoriaId
    $ws.Cells.Item($row, 10).Value = $categoria
    $ws.Cells.Item($row, 11).Value = $variedad
    $ws.Cells.Item($row, 12).Value = $calidad
    $ws.Cells.Item($row, 13).Value = $volumen
    $ws.Cells.Item($row, 14).Value = $precioMin
    $ws.Cells.Item($row, 15).Value = $precioMax
    $ws.Cells.Item($row, 16).Value = $precioProm
    $ws.Cells.Item($row, 17).Value = $unidad
    $ws.Cells.Item($row, 18).Value = $origen
    $ws.Cells.Item($row, 19).Value = $precioKg
    $ws.Cells.Item($row, 20).Value = $kgUnidad
}

Set-FrutillaRow 211 "Especial" 300 19500 20000 19750 2821
Set-FrutillaRow 212 "Primera"  400 15500 16000 15750 2250
Set-FrutillaRow 213 "Segunda"  300 12500 13000 12750 1821
